# LonePair Core part Update
# Refresh the "current" comparison block (rows 27-46) on Sheet1 with the
# latest reference values from the "core" block (rows 5-24), and clear out
# the now-unused per-species weighting table that used to sit alongside it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the Geometric / Internal-Geometric derivative numbers (I:O, rows 27-41)
# down from the refreshed reference rows (5-19). Column L is a blank spacer
# column in both ranges, so including it in the block copy is harmless.
$ws.Range("I27:O41").Value2 = $ws.Range("I5:O19").Value2

# Copy the Lattice Derivative block (D:F, rows 37-39) from rows 15-17.
$ws.Range("D37:F39").Value2 = $ws.Range("D15:F17").Value2

# Copy the Strain Derivative matrix (D:F, rows 44-46) from rows 22-24.
$ws.Range("D44:F46").Value2 = $ws.Range("D22:F24").Value2

# The per-species weighting table (Species/Core-Shell columns Q:U) for the
# lower comparison block is no longer needed now that the two blocks match -
# clear it out (formatting on S:U is retained, matching a Delete-key clear).
$ws.Range("Q48:U63").ClearContents()

# Recalculate so the percentage-difference formulas (rows 49-63, 66-68) pick
# up the now-matching values (they settle to 0 since the blocks are equal).
$excel.Calculate()

# Restore the view: scrolled down two more rows, with K27:K41 selected
# (the block that was just refreshed) instead of the old F42 selection.
$ws.Activate()
$excel.ActiveWindow.TopLeftCell = $ws.Range("A16")
$ws.Range("K27:K41").Select()
